$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update the "Handed back: in sync with en-US" -> "Ready for handoff" status cells
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn.Range("C3").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"

# Update timestamps (stored as text strings, not dates)
$overview.Range("D2").Value = "2016-03-31 07:55:13"
$overview.Range("D3").Value = "2016-03-31 07:55:13"

$zhcn.Range("E2").Value = "2016-03-31 07:55:01"
$zhcn.Range("E3").Value = "2016-03-31 07:55:01"
